# "good accuracy of depth"
# Replace the Distance-Measurement sample rows (A:B, rows 1-38) with the new
# accuracy-of-depth dataset, and append one additional row (row 39) that
# extends the table to A1:B39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple is (row, A-value, B-value) taken from the updated dataset.
$data = @(
    @(1,  3.73, 0.1244032118055556),
    @(2,  3.28, 0.1709526909722222),
    @(3,  4.62, 0.04584418402777778),
    @(4,  3.88, 0.1103515625),
    @(5,  3.13, 0.1886393229166667),
    @(6,  2.83, 0.2251519097222222),
    @(7,  3,    0.2038845486111111),
    @(8,  5.14, 0.006076388888888889),
    @(9,  5.22, 0),
    @(10, 4.79, 0.03255208333333334),
    @(11, 3.41, 0.1571723090277778),
    @(12, 3.25, 0.1745876736111111),
    @(13, 3.75, 0.1227756076388889),
    @(14, 3.29, 0.1706271701388889),
    @(15, 4.02, 0.09727647569444445),
    @(16, 4.33, 0.07042100694444445),
    @(17, 4.27, 0.0749782986111111),
    @(18, 4.32, 0.07118055555555555),
    @(19, 4.3,  0.07275390625),
    @(20, 4.32, 0.0712890625),
    @(21, 2.17, 0.3213975694444444),
    @(22, 2.04, 0.3438042534722222),
    @(23, 4.99, 0.01746961805555556),
    @(24, 5.22, 0),
    @(25, 5.22, 0),
    @(26, 5.22, 0),
    @(27, 4.86, 0.02674696180555556),
    @(28, 5.33, -0.0078125),
    @(29, 5.33, -0.0078125),
    @(30, 5.33, -0.0078125),
    @(31, 5.33, -0.0078125),
    @(32, 4.22, 0.07953559027777778),
    @(33, 4.23, 0.07893880208333333),
    @(34, 4.24, 0.07747395833333333),
    @(35, 4.27, 0.0751953125),
    @(36, 4.28, 0.07416449652777778),
    @(37, 4.3,  0.07242838541666667),
    @(38, 4.3,  0.072265625),
    @(39, 4.18, 0.08295355902777778)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
}
